$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in A1 (45406 -> 45436)
$ws.Range("A1").Value = 45436

# Update the price list values in column D
$ws.Range("D33").Value = 9431.802
$ws.Range("D34").Value = 11160.692
$ws.Range("D35").Value = 13785.194
$ws.Range("D36").Value = 15376.297
$ws.Range("D37").Value = 15901.197

# Re-apply the merged cell ranges (unmerge then merge) so the
# mergeCells collection gets rebuilt/reordered as in the target workbook.
$ws.Range("A1:D1").UnMerge()
$ws.Range("A9:D9").UnMerge()
$ws.Range("B33:C33").UnMerge()
$ws.Range("B34:C34").UnMerge()
$ws.Range("B37:C37").UnMerge()
$ws.Range("B36:C36").UnMerge()
$ws.Range("B32:C32").UnMerge()
$ws.Range("A11:D11").UnMerge()
$ws.Range("A10:D10").UnMerge()
$ws.Range("B35:C35").UnMerge()

$ws.Range("A1:D1").Merge()
$ws.Range("A9:D9").Merge()
$ws.Range("B33:C33").Merge()
$ws.Range("B34:C34").Merge()
$ws.Range("B37:C37").Merge()
$ws.Range("B36:C36").Merge()
$ws.Range("B32:C32").Merge()
$ws.Range("A11:D11").Merge()
$ws.Range("A10:D10").Merge()
$ws.Range("B35:C35").Merge()
